$wb = $excel.ActiveWorkbook

# Update the input probabilities on the "potential_preg_untrt" sheet.
# These values feed downstream formulas on potential_preg_trt,
# potential_preec_untrt, potential_preec_trt and postpreec_preg, which
# will recalculate automatically.
$ws = $wb.Worksheets.Item("potential_preg_untrt")

$ws.Range("C9").Value = 0.05
$ws.Range("C10").Value = 0.02
$ws.Range("C11").Value = 0.02
$ws.Range("C13").Value = 0.005
$ws.Range("C14").Value = 0.004
$ws.Range("C15").Value = 0.004
$ws.Range("C16").Value = 0.004
$ws.Range("C17").Value = 0.004

# Move the active selection/tab from SimParameters to
# potential_preg_untrt, selecting C2:C21.
$ws.Activate()
$ws.Range("C2:C21").Select()
